$d = $word.ActiveDocument

$d.Content.Find.Execute("NTT Data Romania", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NTT Data Romania S.A.", 2)

$d.Content.Find.Execute("Street Constanta 19-21 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "19-21, Constanta Street,", 2)

$d.Content.Find.Execute("Cluj Napoca City, 400158", $true, $false, $false, $false, $false,
                         $true, 1, $false, "400158 Cluj Napoca", 2)
